# "a lot of figures are corrected"
#
# The reflection/moment-of-inertia diagram on slide 1 carried two stray
# delta-e annotation labels ("TextBox 14" -> "Δe.x" and "TextBox 22" -> "Δe.y")
# that duplicated/conflicted with the corrected figure. Remove them.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$namesToRemove = @("TextBox 14", "TextBox 22")

# Collect first (mutating Shapes while iterating forward can skip items),
# then delete.
$shapesToDelete = @()
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($namesToRemove -contains $shp.Name) {
        $shapesToDelete += $shp
    }
}

foreach ($shp in $shapesToDelete) {
    $shp.Delete()
}
